$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: update payment_tokens usd_price / eth_price values
$ws.Range("B2").Value = '[{''id'': 13689077, ''symbol'': ''ETH'', ''address'': ''0x0000000000000000000000000000000000000000'', ''image_url'': ''https://storage.opensea.io/files/6f8e2979d428180222796ff4a33ab929.svg'', ''name'': ''Ether'', ''decimals'': 18, ''eth_price'': 1.0, ''usd_price'': 2064.03}, {''id'': 12182941, ''symbol'': ''DAI'', ''address'': ''0x6b175474e89094c44da98b954eedeac495271d0f'', ''image_url'': ''https://storage.opensea.io/files/8ef8fb3fe707f693e57cdbfea130c24c.svg'', ''name'': ''Dai Stablecoin'', ''decimals'': 18, ''eth_price'': 0.00048448908203853623, ''usd_price'': 1.0}, {''id'': 4645681, ''symbol'': ''WETH'', ''address'': ''0xc02aaa39b223fe8d0a0e5c4f27ead9083c756cc2'', ''image_url'': ''https://storage.opensea.io/files/accae6b6fb3888cbff27a013729c22dc.svg'', ''name'': ''Wrapped Ether'', ''decimals'': 18, ''eth_price'': 1.0, ''usd_price'': 2064.03}, {''id'': 4403908, ''symbol'': ''USDC'', ''address'': ''0xa0b86991c6218b36c1d19d4a2e9eb0ce3606eb48'', ''image_url'': ''https://storage.opensea.io/files/749015f009a66abcb3bbb3502ae2f1ce.svg'', ''name'': ''USD Coin'', ''decimals'': 6, ''eth_price'': 0.0004842700929734548, ''usd_price'': 0.999548}]'

# Row 2 stats.* numeric updates (AG2:BA2)
$ws.Range("AG2").Value = 117.664697575321
$ws.Range("AH2").Value = 0.6284829264734667
$ws.Range("AI2").Value = 317
$ws.Range("AJ2").Value = 0.3711820112786152
$ws.Range("AK2").Value = 2104.73656844445
$ws.Range("AM2").Value = 4341
$ws.Range("AN2").Value = 0.4848506262253973
$ws.Range("AO2").Value = 2104.73656844445
$ws.Range("AQ2").Value = 4341
$ws.Range("AR2").Value = 0.4848506262253973
$ws.Range("AS2").Value = 2104.73656844445
$ws.Range("AT2").Value = 4343
$ws.Range("AW2").Value = 1912
$ws.Range("AX2").Value = 0.4846273470975017
$ws.Range("AZ2").Value = 4848.506262253974
$ws.Range("BA2").Value = 0.269
